$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 6 data: Dead volume label, value, and unit
$ws.Range("A6").Value = "Dead volume"
$ws.Range("B6").Value = 0.75
$ws.Range("C6").Value = "mL"

# Apply style to B6: fill (yellow) + right alignment, matching B2:B5 style plus right align
$ws.Range("B6").Interior.Color = $ws.Range("B5").Interior.Color
$ws.Range("B6").HorizontalAlignment = -4152  # xlRight

# Update selection to match the new active cell
$ws.Range("B6").Select()
